$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 3.8
$ws.Range("L3").Value = 4.5
$ws.Range("W3").Value = 6
$ws.Range("AH3").Value = 17
$ws.Range("AV3").Value = 67
$ws.Range("AY3").Value = 34
